$d = $word.ActiveDocument

# --- Header table: Enquiry Date ---
# Table 1, Row 2 ("Enquiry" | "a" | "Date" | "12/20/2019, 04:40 PM")
$t1 = $d.Tables.Item(1)
$t1.Cell(2, 4).Range.Text = "12/25/2019, 10:49 AM"

# --- COOLING WATER CIRCUIT table (Table 4) ---
$t4 = $d.Tables.Item(4)
# Row 2: "1." | "Cooling water flow" | "m3/hr" | "125" -> "114"
$t4.Cell(2, 4).Range.Text = "114"
# Row 4: "3." | "Cooling water outlet temperature" | "C" | "36.6" -> "37.1"
$t4.Cell(4, 4).Range.Text = "37.1"
# Row 7: "6." | "Cooling water circuit pressure loss" | "mLC" | "2.6" -> "2.2"
$t4.Cell(7, 4).Range.Text = "2.2"

# --- Steam Circuit table (Table 5) ---
$t5 = $d.Tables.Item(5)
# Row 3: "2." | "Steam Consumption(+/-3%)" | "kg/hr" | "398.1" -> "400.5"
$t5.Cell(3, 4).Range.Text = "400.5"
